$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sem" (sheet1): add a new "Score20" column (T) with per-row data.
# ---------------------------------------------------------------------------
$sem = $wb.Worksheets.Item("Sem")
$sem.Activate()

# Header
$sem.Range("T1").Value = "Score20"

# Row data (PS_No. rows 2..16)
$sem.Range("T2").Value = 91
$sem.Range("T3").Value = 43
$sem.Range("T4").Value = 39
$sem.Range("T5").Value = 39
$sem.Range("T6").Value = 65
$sem.Range("T7").Value = 75
$sem.Range("T8").Value = 65
$sem.Range("T9").Value = 39
$sem.Range("T10").Value = 43
$sem.Range("T11").Value = 55
$sem.Range("T12").Value = 93
$sem.Range("T13").Value = 99
$sem.Range("T14").Value = 93
$sem.Range("T15").Value = 39
$sem.Range("T16").Value = 76

# Scroll the window so column E is the left-most visible column, then leave
# the just-entered T12:T16 block selected (mirrors the editor's final state
# after typing in the new column).
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$sem.Range("T12:T16").Select()

# ---------------------------------------------------------------------------
# Remaining (empty) sheets: the editor had previously selected T12:T16 on
# "Sem" and that multi-area selection carried over as it tabbed through the
# other sheets, each landing back on their own original active cell.
# ---------------------------------------------------------------------------
$hobbies = $wb.Worksheets.Item("Hobbies")
$hobbies.Activate()
$u1 = $excel.Union($hobbies.Range("T12:T16"), $hobbies.Range("A1"))
$u1.Select()

$cities = $wb.Worksheets.Item("Cities")
$cities.Activate()
$u2 = $excel.Union($cities.Range("T12:T16"), $cities.Range("A1"))
$u2.Select()

$pl = $wb.Worksheets.Item("PL")
$pl.Activate()
$u3 = $excel.Union($pl.Range("T12:T16"), $pl.Range("A1"))
$u3.Select()

$domain = $wb.Worksheets.Item("Domain")
$domain.Activate()
$u4 = $excel.Union($domain.Range("T12:T16"), $domain.Range("A4"))
$u4.Select()

# End back on the "Sem" sheet, which is the tab that was active when saved.
$sem.Activate()
$sem.Range("T12:T16").Select()
